# Add 2022-Q3 data:
#  - Insert a new "2022-Q3" summary row at the top of the "总计" sheet,
#    pushing the existing quarters down by one row.
#  - Insert a new "2022-Q3" worksheet (positioned right after "总计",
#    before "2022-Q2") holding the per-fund holdings detail for the
#    new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: shift the existing quarters down one row
#    and write the new 2022-Q3 figures into row 2. Values are written
#    as literals (rather than copied cell-to-cell) to avoid any
#    COM Variant round-trip float noise.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give row 9 (brand new row) the same formatting (borders etc.) as row 8
# before putting values into it, so style indices line up with the rest
# of the table.
$summary.Range("A8:D8").Copy()
$summary.Range("A9:D9").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q3", 4, 1.26),
    @(1, "2022-Q2", 4, 1.85),
    @(2, "2022-Q1", 6, 1.71),
    @(3, "2021-Q4", 4, 1.22),
    @(4, "2021-Q3", 5, 1.77),
    @(5, "2021-Q2", 2, 1.82),
    @(6, "2021-Q1", 3, 1.05),
    @(7, "2020-Q4", 3, 1.96)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Range("A$r").Value = $row[0]
    $summary.Range("B$r").Value = $row[1]
    $summary.Range("C$r").Value = $row[2]
    $summary.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (so it keeps
#    identical headers/formatting) positioned right before it, rename it,
#    then overwrite the data rows with the new quarter's numbers.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The source columns are stored as text (fund code/name/size/position
# numbers all came in as strings), so force text formatting before
# writing the values to avoid Excel auto-converting them to numbers.
$q3.Range("B2:G5").NumberFormat = "@"

$q3.Range("B2").Value = "002207"
$q3.Range("C2").Value = "前海开源金银珠宝主题精选混合C"
$q3.Range("D2").Value = "6.72"
$q3.Range("E2").Value = "90.85"
$q3.Range("F2").Value = "7.85"
$q3.Range("G2").Value = "0.5275"
$q3.Range("H2").Value = 7

$q3.Range("B3").Value = "001302"
$q3.Range("C3").Value = "前海开源金银珠宝主题精选混合A"
$q3.Range("D3").Value = "3.99"
$q3.Range("E3").Value = "90.85"
$q3.Range("F3").Value = "7.85"
$q3.Range("G3").Value = "0.3132"
$q3.Range("H3").Value = 7

$q3.Range("B4").Value = "003304"
$q3.Range("C4").Value = "前海开源沪港深核心资源灵活配置混合A"
$q3.Range("D4").Value = "3.45"
$q3.Range("E4").Value = "90.59"
$q3.Range("F4").Value = "7.80"
$q3.Range("G4").Value = "0.2691"
$q3.Range("H4").Value = 6

$q3.Range("B5").Value = "003305"
$q3.Range("C5").Value = "前海开源沪港深核心资源灵活配置混合C"
$q3.Range("D5").Value = "1.89"
$q3.Range("E5").Value = "90.59"
$q3.Range("F5").Value = "7.80"
$q3.Range("G5").Value = "0.1474"
$q3.Range("H5").Value = 6

# Copying a sheet makes the copy the active one; restore the original
# active sheet ("总计") so the workbook-level active-tab pointer is
# unchanged, same as before the edit.
$summary.Activate()
